$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the description for the January 2017 county update row (row 3, col C)
$ws.Range("C3").Value = "county updates 01/2017(released 2017-03-01)"

# Add the new row for the February 2017 county update (row 4)
$ws.Range("A4").Value = "2017年2月中华人民共和国县以上行政区划代码.csv"
$ws.Range("B4").Value = "http://www.mca.gov.cn/article/sj/tjbz/a/2017/0327/2017%E5%B9%B42%E6%9C%88%E4%B8%AD%E5%8D%8E%E4%BA%BA%E6%B0%91%E5%85%B1%E5%92%8C%E5%9B%BD%E5%8E%BF%E4%BB%A5%E4%B8%8A%E8%A1%8C%E6%94%BF%E5%8C%BA%E5%88%92%E4%BB%A3%E7%A0%81.html"
$ws.Range("C4").Value = "county updates 02/2017(released 2017-03-27)"

# Match the selection shown in the saved file
$ws.Range("C2").Select()

# The saved workbook also carries a page-setup definition (portrait orientation)
$ws.PageSetup.Orientation = 1
